$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values to reflect the updated election results
$ws.Range("H2").Value = 56
$ws.Range("I2").Value = 178
$ws.Range("J2").Value = 722
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 198
$ws.Range("M2").Value = 13
$ws.Range("N2").Value = 125
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 4
$ws.Range("R2").Value = 9
$ws.Range("S2").Value = 85
$ws.Range("T2").Value = 137
$ws.Range("U2").Value = 5
$ws.Range("V2").Value = 1114
$ws.Range("X2").Value = 1107
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 19
$ws.Range("AA2").Value = 8

$wb.Save()
